$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.902.27"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.631.95"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "1.865.34"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "1.630.34"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").Value = "29.902.92"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "1.422.94"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.89%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.556"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.835"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0502"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.772.39"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("E51").Value = "  +3.79%  "
